# Refresh the "cryptos" price/volume snapshot (D2:E51) with the latest
# scraped values. Price values that look numeric are entered with a
# leading apostrophe (forces text, like the source scraper's raw
# inline-string cells) and then restyled back to "Normal" so no stray
# Text number-format/quote-prefix style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.239.72'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.863.25'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'0.7123"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '
$ws.Range("D6").Value = "'240.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").Value = "'0.07710"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("D10").Value = "'24.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("D11").Value = "'0.08352"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '1.867.17'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").Value = "'5.202"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").Value = "'0.7131"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'91.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '29.243.65'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = "'5.951"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = "'242.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").Value = "'0.000007837"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").Value = '2.122.32'
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = "'7.844"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.15%  '
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = "'0.1594"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = "'163.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").Value = "'8.890"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").Value = "'1.498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = "'4.409"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").Value = "'4.244"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("D33").Value = "'0.05139"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").Value = "'0.8082"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.53%  '
$ws.Range("D35").Value = "'1.933"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = "'1.168"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = "'0.01852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").Value = '1.176.51'
$ws.Range("E40").Value = '  -6.79%  '
$ws.Range("D41").Value = "'6.207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.56%  '
$ws.Range("D42").Value = "'0.8936"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'72.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = "'102.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").Value = '2.019.84'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").Value = "'0.5193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.51%  '
$ws.Range("D48").Value = "'1.791"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = "'9.283"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = "'0.9998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.20%  '
